# Natmi LR-pairs results for Spon2-Itgam were recomputed ("Natmi following Dr Hou advice").
# The sending/target cluster set grew from {FAPs, sCs} to {ECs, FAPs, sCs}, so the 2x2
# combination table (rows 2-5) becomes a 3x2 combination table (rows 2-7), and every
# numeric statistic (ligand/receptor expression, specificity, edge weights, ...) is refreshed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A=Sending cluster, B=Ligand symbol, C=Receptor symbol, D=Target cluster,
# E=Ligand-expressing cells, F=Ligand detection rate, G=Ligand average expression value,
# H=Ligand total expression value, I=Ligand derived specificity (avg), J=Ligand derived specificity (total),
# K=Receptor-expressing cells, L=Receptor detection rate, M=Receptor average expression value,
# N=Receptor total expression value, O=Receptor derived specificity (avg), P=Receptor derived specificity (total),
# Q=Edge average expression weight, R=Edge total expression weight,
# S=Edge average expression derived specificity, T=Edge total expression derived specificity
$data = @(
    @("ECs", "Spon2", "Itgam", "M2", 2, 0.6666666666666666, 0.436226, 1.308678, 0.02612337021374315, 0.02612337021374314, 3, 1, 45.931316, 137.793948, 0.9874217014725413, 0.9874217014725412, 20.036434253416, 180.327908280744, 0.02579478266465136, 0.02579478266465136),
    @("ECs", "Spon2", "Itgam", "sCs", 2, 0.6666666666666666, 0.436226, 1.308678, 0.02612337021374315, 0.02612337021374314, 3, 1, 0.5850973333333334, 1.755292, 0.01257829852745884, 0.01257829852745884, 0.2552346693306667, 2.297112023976, 0.0003285875490917875, 0.0003285875490917875),
    @("FAPs", "Spon2", "Itgam", "M2", 3, 1, 15.76143266666667, 47.284298, 0.9438725354525366, 0.9438725354525365, 3, 1, 45.931316, 137.793948, 0.9874217014725413, 0.9874217014725412, 723.9433444253893, 6515.490099828504, 0.9320002249297452, 0.932000224929745),
    @("FAPs", "Spon2", "Itgam", "sCs", 3, 1, 15.76143266666667, 47.284298, 0.9438725354525366, 0.9438725354525365, 3, 1, 0.5850973333333334, 1.755292, 0.01257829852745884, 0.01257829852745884, 9.221972222779556, 82.997750005016, 0.01187231052279148, 0.01187231052279148),
    @("sCs", "Spon2", "Itgam", "M2", 3, 1, 0.5010290000000001, 1.503087, 0.0300040943337204, 0.0300040943337204, 3, 1, 45.931316, 137.793948, 0.9874217014725413, 0.9874217014725412, 23.012921324164, 207.116291917476, 0.02962669387814484, 0.02962669387814483),
    @("sCs", "Spon2", "Itgam", "sCs", 3, 1, 0.5010290000000001, 1.503087, 0.0300040943337204, 0.0300040943337204, 3, 1, 0.5850973333333334, 1.755292, 0.01257829852745884, 0.01257829852745884, 0.2931507318226667, 2.638356586404, 0.0003774004555755714, 0.0003774004555755714)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

Write-Output "Done. UsedRange:"
Write-Output $ws.UsedRange.Address()

